$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("formData")

# New test case row (row 5) mirroring the style used by the existing data rows (2-4)
$ws.Range("A5").Value = "TC04"
$ws.Range("F5").Value = "Blackberry"
$ws.Range("G5").Value = "India"
$ws.Range("H5").Value = "ind"

$src = $ws.Range("A4")
$src.Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("H5").PasteSpecial(-4122)
